$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102, shifting existing rows 102-122 down to 103-123.
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with the new weekly record.
$ws.Cells.Item(102, 1).Value = 1
$ws.Cells.Item(102, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(102, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(102, 4).Value = 45131
$ws.Cells.Item(102, 4).NumberFormat = $ws.Cells.Item(103, 4).NumberFormat
$ws.Cells.Item(102, 5).Value = 15
$ws.Cells.Item(102, 6).Value = 100112040
$ws.Cells.Item(102, 7).Value = "Cilantro"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 300
$ws.Cells.Item(102, 11).Value = 1800
$ws.Cells.Item(102, 12).Value = 2000
$ws.Cells.Item(102, 13).Value = 1900
$ws.Cells.Item(102, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(102, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(102, 16).Value = 950
$ws.Cells.Item(102, 17).Value = 2
$ws.Cells.Item(102, 18).Value = "Hortaliza"
